$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.9736842105263158
$ws.Range("D3").Value = 1
$ws.Range("H3").Value = 0.8886443661971831
$ws.Range("I3").Value = 0.02697572826874109
$ws.Range("J3").Value = 0.868421052631579
$ws.Range("K3").Value = 141.0526315789474
$ws.Range("Q3").Value = 56
$ws.Range("R3").Value = 78
$ws.Range("S3").Value = 110
$ws.Range("T3").Value = 122
$ws.Range("U3").Value = 151
$ws.Range("V3").Value = 4450
$ws.Range("W3").Value = 4428
$ws.Range("X3").Value = 4396
$ws.Range("Y3").Value = 4384
$ws.Range("Z3").Value = 4355
$ws.Range("AF3").Value = 0.987572
$ws.Range("AG3").Value = 0.98269
$ws.Range("AH3").Value = 0.975588
$ws.Range("AI3").Value = 0.972925
$ws.Range("AJ3").Value = 0.966489
